# Generate Report for Handoff
# Adds two new rows (a .png asset handoff and the final .md handoff) to each of the
# three worksheets (Overview, zh-cn, de-de), refreshes the existing row's values, and
# rebuilds the hyperlinks collections to point at the new files.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "1764e6c3-3e2e-41da-9573-703a6d3eed43.png"
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"
$ws1.Range("D2").Value = "2016-03-24 13:20:41"

$ws1.Range("A3").Value = "7f6da85b-2e8a-4057-a82c-90145b000f2a.png"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"
$ws1.Range("D3").Value = "2016-03-24 13:20:41"
$ws1.Range("A3").Style = "HyperLink"
$ws1.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws1.Range("A4").Value = "b6ada119-9f9e-4991-a1db-da1d89996ff8.md"
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"
$ws1.Range("D4").Value = "2016-03-24 13:20:41"
$ws1.Range("A4").Style = "HyperLink"
$ws1.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/bdb3449da203ace61162dcccd403b25786529026/e2e/1764e6c3-3e2e-41da-9573-703a6d3eed43.png", "", "", "1764e6c3-3e2e-41da-9573-703a6d3eed43.png")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/bdb3449da203ace61162dcccd403b25786529026/e2e/7f6da85b-2e8a-4057-a82c-90145b000f2a.png", "", "", "7f6da85b-2e8a-4057-a82c-90145b000f2a.png")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/bdb3449da203ace61162dcccd403b25786529026/e2e/b6ada119-9f9e-4991-a1db-da1d89996ff8.md", "", "", "b6ada119-9f9e-4991-a1db-da1d89996ff8.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "1764e6c3-3e2e-41da-9573-703a6d3eed43.png"
$ws2.Range("B2").Value = ".png"
$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("D2").Value = "c22b45e6bcc0874f7f84130fdfd6a88f79c81c1d.png"
$ws2.Range("E2").Value = "2016-03-24 13:20:37"
$ws2.Range("H2").Value = "0001-01-01 00:00:00"
$ws2.Range("J2").Value = "IsDependency"
$ws2.Range("K2").Value = "e2e\b6ada119-9f9e-4991-a1db-da1d89996ff8.md"

$ws2.Range("A3").Value = "7f6da85b-2e8a-4057-a82c-90145b000f2a.png"
$ws2.Range("B3").Value = ".png"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("D3").Value = "d170302f4b58da2541c93573c0d30bc86ff4edc6.png"
$ws2.Range("E3").Value = "2016-03-24 13:20:37"
$ws2.Range("H3").Value = "0001-01-01 00:00:00"
$ws2.Range("J3").Value = "IsDependency"
$ws2.Range("K3").Value = "e2e\b6ada119-9f9e-4991-a1db-da1d89996ff8.md"
$ws2.Range("A3").Style = "HyperLink"
$ws2.Range("D3").Style = "HyperLink"
$ws2.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws2.Range("A4").Value = "b6ada119-9f9e-4991-a1db-da1d89996ff8.md"
$ws2.Range("B4").Value = ".md"
$ws2.Range("C4").Value = "Ready for handoff"
$ws2.Range("D4").Value = "b6ada119-9f9e-4991-a1db-da1d89996ff8.d3b970b1bc158ed7223b89b28c3db9fc2b229eca.zh-cn.xlf"
$ws2.Range("E4").Value = "2016-03-24 13:20:37"
$ws2.Range("H4").Value = "0001-01-01 00:00:00"
$ws2.Range("J4").Value = "Include"
$ws2.Range("A4").Style = "HyperLink"
$ws2.Range("D4").Style = "HyperLink"
$ws2.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/bdb3449da203ace61162dcccd403b25786529026/e2e/1764e6c3-3e2e-41da-9573-703a6d3eed43.png", "", "", "1764e6c3-3e2e-41da-9573-703a6d3eed43.png")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/80f35451704f938735b3cec31209ac9aec4f5258/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c22b45e6bcc0874f7f84130fdfd6a88f79c81c1d.png", "", "", "c22b45e6bcc0874f7f84130fdfd6a88f79c81c1d.png")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/bdb3449da203ace61162dcccd403b25786529026/e2e/7f6da85b-2e8a-4057-a82c-90145b000f2a.png", "", "", "7f6da85b-2e8a-4057-a82c-90145b000f2a.png")
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/80f35451704f938735b3cec31209ac9aec4f5258/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d170302f4b58da2541c93573c0d30bc86ff4edc6.png", "", "", "d170302f4b58da2541c93573c0d30bc86ff4edc6.png")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/bdb3449da203ace61162dcccd403b25786529026/e2e/b6ada119-9f9e-4991-a1db-da1d89996ff8.md", "", "", "b6ada119-9f9e-4991-a1db-da1d89996ff8.md")
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/80f35451704f938735b3cec31209ac9aec4f5258/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b6ada119-9f9e-4991-a1db-da1d89996ff8.d3b970b1bc158ed7223b89b28c3db9fc2b229eca.zh-cn.xlf", "", "", "b6ada119-9f9e-4991-a1db-da1d89996ff8.d3b970b1bc158ed7223b89b28c3db9fc2b229eca.zh-cn.xlf")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "1764e6c3-3e2e-41da-9573-703a6d3eed43.png"
$ws3.Range("B2").Value = ".png"
$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("D2").Value = "c22b45e6bcc0874f7f84130fdfd6a88f79c81c1d.png"
$ws3.Range("E2").Value = "2016-03-24 13:20:41"
$ws3.Range("H2").Value = "0001-01-01 00:00:00"
$ws3.Range("J2").Value = "IsDependency"
$ws3.Range("K2").Value = "e2e\b6ada119-9f9e-4991-a1db-da1d89996ff8.md"

$ws3.Range("A3").Value = "7f6da85b-2e8a-4057-a82c-90145b000f2a.png"
$ws3.Range("B3").Value = ".png"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("D3").Value = "d170302f4b58da2541c93573c0d30bc86ff4edc6.png"
$ws3.Range("E3").Value = "2016-03-24 13:20:41"
$ws3.Range("H3").Value = "0001-01-01 00:00:00"
$ws3.Range("J3").Value = "IsDependency"
$ws3.Range("K3").Value = "e2e\b6ada119-9f9e-4991-a1db-da1d89996ff8.md"
$ws3.Range("A3").Style = "HyperLink"
$ws3.Range("D3").Style = "HyperLink"

$ws3.Range("A4").Value = "b6ada119-9f9e-4991-a1db-da1d89996ff8.md"
$ws3.Range("B4").Value = ".md"
$ws3.Range("C4").Value = "Ready for handoff"
$ws3.Range("D4").Value = "b6ada119-9f9e-4991-a1db-da1d89996ff8.d3b970b1bc158ed7223b89b28c3db9fc2b229eca.de-de.xlf"
$ws3.Range("E4").Value = "2016-03-24 13:20:41"
$ws3.Range("H4").Value = "0001-01-01 00:00:00"
$ws3.Range("J4").Value = "Include"
$ws3.Range("A4").Style = "HyperLink"
$ws3.Range("D4").Style = "HyperLink"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/bdb3449da203ace61162dcccd403b25786529026/e2e/1764e6c3-3e2e-41da-9573-703a6d3eed43.png", "", "", "1764e6c3-3e2e-41da-9573-703a6d3eed43.png")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/60f13a73b9b8e07da096745d33043af68135c935/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c22b45e6bcc0874f7f84130fdfd6a88f79c81c1d.png", "", "", "c22b45e6bcc0874f7f84130fdfd6a88f79c81c1d.png")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/bdb3449da203ace61162dcccd403b25786529026/e2e/7f6da85b-2e8a-4057-a82c-90145b000f2a.png", "", "", "7f6da85b-2e8a-4057-a82c-90145b000f2a.png")
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/60f13a73b9b8e07da096745d33043af68135c935/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d170302f4b58da2541c93573c0d30bc86ff4edc6.png", "", "", "d170302f4b58da2541c93573c0d30bc86ff4edc6.png")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/bdb3449da203ace61162dcccd403b25786529026/e2e/b6ada119-9f9e-4991-a1db-da1d89996ff8.md", "", "", "b6ada119-9f9e-4991-a1db-da1d89996ff8.md")
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/60f13a73b9b8e07da096745d33043af68135c935/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b6ada119-9f9e-4991-a1db-da1d89996ff8.d3b970b1bc158ed7223b89b28c3db9fc2b229eca.de-de.xlf", "", "", "b6ada119-9f9e-4991-a1db-da1d89996ff8.d3b970b1bc158ed7223b89b28c3db9fc2b229eca.de-de.xlf")

Write-Host "Report generated for handoff: added rows for 1764e6c3-3e2e-41da-9573-703a6d3eed43.png, 7f6da85b-2e8a-4057-a82c-90145b000f2a.png, b6ada119-9f9e-4991-a1db-da1d89996ff8.md"
